$wb = $excel.ActiveWorkbook

# --- Grab sheet references before renaming ---
$wsSearch  = $wb.Worksheets.Item("search")
$wsTickets = $wb.Worksheets.Item("tickets")
$wsData    = $wb.Worksheets.Item("data")

# --- Rename sheets: tickets -> cabin, data -> user ---
$wsTickets.Name = "cabin"
$wsData.Name = "user"

# --- Update "cabin" (ex "tickets") sheet data ---
# Row 3: was (1,0,1) -> becomes (2,1,1)
$wsTickets.Cells.Item(3, 1).Value = 2
$wsTickets.Cells.Item(3, 2).Value = 1
$wsTickets.Cells.Item(3, 3).Value = 1

# New row 4 (old row 3 values shifted down): (1,0,1)
$wsTickets.Cells.Item(4, 1).Value = 1
$wsTickets.Cells.Item(4, 2).Value = 0
$wsTickets.Cells.Item(4, 3).Value = 1

# New row 5 (old row 4 values shifted down): (1,5,5)
$wsTickets.Cells.Item(5, 1).Value = 1
$wsTickets.Cells.Item(5, 2).Value = 5
$wsTickets.Cells.Item(5, 3).Value = 5

# New formatted (empty) cells E3 and I17 - underline font style
$wsTickets.Range("E3").Font.Underline = $true
$wsTickets.Range("I17").Font.Underline = $true

# Page setup: portrait orientation
$wsTickets.PageSetup.Orientation = 1

# --- Update selections (without permanently changing the active sheet) ---
$excel.Goto($wsTickets.Range("E3"))
$excel.Goto($wsData.Range("F3"))

Write-Host ("Sheets: " + (($wb.Worksheets | ForEach-Object { $_.Name }) -join ", "))
